# Weekly update: a new price-report entry is added for the "Vega Modelo de
# Temuco - Espárragos" sheet. It becomes the new row 36, pushing the
# existing rows 36..123 down to 37..124 (the sheet is a daily/weekly log
# that keeps growing by one row at the top of the insertion point).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 36; everything below (old rows 36..123) shifts
# down to 37..124, and the used range grows to A1:R124 automatically.
$ws.Rows.Item(36).Insert()

# The new row reuses the same classification fields (market, region,
# category, variety, quality, unit, origin, kg/unit, classification) as
# the entry that is now sitting at row 37 (the old row 36) - only the
# date/volume/price columns differ for the new entry. Copy that row down
# into the freshly inserted one as a starting point.
$ws.Rows.Item(37).Copy()
$ws.Rows.Item(36).PasteSpecial()
$excel.CutCopyMode = 0

# Now overwrite the new entry's date, volume and price columns.
$ws.Cells.Item(36, 4).Value  = 45260   # D36 - Fecha
$ws.Cells.Item(36, 10).Value = 300     # J36 - Volumen
$ws.Cells.Item(36, 11).Value = 1800    # K36 - Precio minimo
$ws.Cells.Item(36, 12).Value = 1800    # L36 - Precio maximo
$ws.Cells.Item(36, 13).Value = 1800    # M36 - Precio promedio ponderado
$ws.Cells.Item(36, 16).Value = 1800    # P36 - Precio $/Kg
